$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Testcase "001_Login_002_LoginWithoutEmailAndPassword" (row 12):
# mark the login-button check step, and flip its "Check" action off (<NOP>)
# now that a second testcase continues further down.
$ws.Range("J12").Value = "X"
$ws.Range("L12").Value = "<NOP>"

# Start of the new, second login testcase (row 17): set the "Funktioniert
# noch nicht" (not working yet) step and click the login button.
$ws.Range("A17").Value = "Funktioniert noch nicht"
$ws.Range("B17").Value = "<SET>"
$ws.Range("L17").Value = "butEinloggen"

# Move the selection/view to just past the new last row, matching where
# the author left the cursor after typing the new row.
$ws.Range("L18").Select()
